# Edit script for LOQ4201.xlsx
# - Inserts a new row at position 13 (the previously-missing row that holds
#   the 'Docentes responsaveis:' content next to its label), shifting the
#   rows below it (old 13-21) down to 14-22.
# - Fills in the syllabus/program/bibliography content that was missing or
#   mis-referenced in several B/C cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(13).Insert()

# The inserted row 13 has no data in columns B/C (it previously held only the
# 'Programa resumido:' label in column A, which has now shifted to row 14),
# and the stray, empty, styled A13 cell left behind by Insert() needs to go
# so the row matches the target (no A13 cell at all).
$ws.Range("A13").Clear()

# B13/C13 never had formatting of their own (row 13 was a label-only row
# before), so copy the wrap/vertical-top formatting used throughout column B
# and C from a fully-populated row before setting the values.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null

$ws.Range("B10").Value = 'Apresentar a engenharia de produção e suas principais áreas de atuação, bem como introduzir o aluno num ambiente de engenharia.'
$ws.Range("C10").Value = 'Apresentar a engenharia de produção e suas principais áreas de atuação, bem como introduzir o aluno num ambiente de engenharia.'
$ws.Range("B13").Value = '11079086 - Herlandí de Souza Andrade'
$ws.Range("C13").Value = '11079086 - Herlandí de Souza Andrade'
$ws.Range("B14").Value = 'A engenharia de Produção. As grandes áreas da engenharia de produção. Projeto de Engenharia.'
$ws.Range("C14").Value = 'A engenharia de Produção. As grandes áreas da engenharia de produção. Projeto de Engenharia.'
$ws.Range("B16").Value = '1 - A ENGENHARIA DE PRODUÇÃO: O campo de trabalho do engenheiro de produção..As áreas de atuação da engenharia de produção. O currículo do curso de engenharia de produção na EEL/USP' + [char]10 + '2 - AS GRANDES AREAS DA ENGENHARIA DE PRODUÇÃO: Noções básicas de Planejamento e controle da produção, pesquisa operacional. Gestão e controle de qualidade. Projeto do produto. Projeto da fábrica. Projeto e estudo de métodos de trabalho. Engenharia da Sustentabilidade, Engenharia Econômica, Gestão de Tecnologia de Informação.' + [char]10 + '3  PROJETO DE ENGENHARIA: Noções básicas de projetos de engenharia.'
$ws.Range("C16").Value = '1 - A ENGENHARIA DE PRODUÇÃO: O campo de trabalho do engenheiro de produção..As áreas de atuação da engenharia de produção. O currículo do curso de engenharia de produção na EEL/USP' + [char]10 + '2 - AS GRANDES AREAS DA ENGENHARIA DE PRODUÇÃO: Noções básicas de Planejamento e controle da produção, pesquisa operacional. Gestão e controle de qualidade. Projeto do produto. Projeto da fábrica. Projeto e estudo de métodos de trabalho. Engenharia da Sustentabilidade, Engenharia Econômica, Gestão de Tecnologia de Informação.' + [char]10 + '3  PROJETO DE ENGENHARIA: Noções básicas de projetos de engenharia.'
$ws.Range("B19").Value = 'Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras'
$ws.Range("C19").Value = 'Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras'
$ws.Range("B20").Value = 'Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas.'
$ws.Range("C20").Value = 'Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas.'
$ws.Range("B21").Value = 'NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.'
$ws.Range("C21").Value = 'NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.'
$ws.Range("B22").Value = 'BAZZO, W.A. Introdução à Engenharia. Florianópolis: Editora da UFSC, 1998.' + [char]10 + 'SLACK, N. et. al. Administração da Produção. São Paulo, Atlas, 1999. ' + [char]10 + 'BATALHA, M.O. et al , Introdução à Engenharia de Produção, Rio de Janeiro, Elsevier, 2008.'
$ws.Range("C22").Value = 'BAZZO, W.A. Introdução à Engenharia. Florianópolis: Editora da UFSC, 1998.' + [char]10 + 'SLACK, N. et. al. Administração da Produção. São Paulo, Atlas, 1999. ' + [char]10 + 'BATALHA, M.O. et al , Introdução à Engenharia de Produção, Rio de Janeiro, Elsevier, 2008.'
